# Apply the crypto price/volume refresh produced by the GitHub Actions job.
# Rows 15/16 (Polygon <-> WrappedliquidstakedEther2.0) also swap position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.307.74'
$ws.Range('E2').Value = '  -4.64%  '
$ws.Range('D3').Value = '2.242.95'
$ws.Range('E3').Value = '  -5.68%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.43'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.587'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -8.40%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.566'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -8.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.18'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -9.23%  '
$ws.Range('E11').Value = '  -2.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0830'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -9.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.74'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -9.60%  '
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.583.22'
$ws.Range('E15').Value = '  -5.70%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.869'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -11.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.45'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -6.93%  '
$ws.Range('D18').Value = '2.248.46'
$ws.Range('E18').Value = '  -5.11%  '
$ws.Range('D19').Value = '43.248.46'
$ws.Range('E19').Value = '  -4.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.51'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -10.23%  '
$ws.Range('E21').Value = '  -8.95%  '
$ws.Range('E22').Value = '  -11.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.63'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -10.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.20'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -12.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '238.20'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -8.88%  '
$ws.Range('E26').Value = '  -8.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.09'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.08'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -10.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.25'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.15%  '
$ws.Range('E31').Value = '  -15.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '35.76'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.51'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -8.48%  '
$ws.Range('E34').Value = '  -8.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.12'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -7.76%  '
$ws.Range('E36').Value = '  -4.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.19'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +8.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.97'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.69%  '
$ws.Range('E39').Value = '  -7.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.46'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.78%  '
$ws.Range('E41').Value = '  -11.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.72'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -9.03%  '
$ws.Range('E43').Value = '  -8.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.22'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').Value = '1.783.36'
$ws.Range('E46').Value = '  -2.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.89'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -11.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.206'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -10.06%  '
$ws.Range('E49').Value = '  -11.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '76.31'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -9.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.43'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -15.81%  '
